# Cpu - Processed: add Mean increase / Median increase columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the _xlchart defined names so that v1.0-v1.4 point at the
#        ranges that used to be called v1.2/v1.3/v1.4/v1.10/v1.11, and drop
#        the now-unused v1.10 / v1.11 duplicates.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlchart.v1.10") { $n.Delete() }
}
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlchart.v1.11") { $n.Delete() }
}
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlchart.v1.0") { $n.RefersTo = "=Blad1!`$A`$66:`$A`$95" }
    elseif ($n.Name -eq "_xlchart.v1.1") { $n.RefersTo = "=Blad1!`$B`$65" }
    elseif ($n.Name -eq "_xlchart.v1.2") { $n.RefersTo = "=Blad1!`$B`$66:`$B`$95" }
    elseif ($n.Name -eq "_xlchart.v1.3") { $n.RefersTo = "=Blad1!`$A`$98:`$A`$187" }
    elseif ($n.Name -eq "_xlchart.v1.4") { $n.RefersTo = "=Blad1!`$B`$98:`$B`$187" }
}

# --- 2. Add "Mean increase" / "Median increase" header + computed cells for
#        each of the three frequency blocks (Low=17-32, Medium=49-64,
#        High=81-96) plus the summary block (All, 113-128).

function Add-IncreaseBlock($headerRow, $formulaRow, $meanRef, $medianRef, $dFormula, $fFormula) {
    $ws.Cells.Item($headerRow, 4).Value = "Mean increase"
    $ws.Cells.Item($headerRow, 4).Font.Bold = $true
    $ws.Cells.Item($headerRow, 6).Value = "Median increase"
    $ws.Cells.Item($headerRow, 6).Font.Bold = $true

    $dCell = $ws.Cells.Item($formulaRow, 4)
    $dCell.Formula = $dFormula
    $dCell.Style = "Standaard"

    $fCell = $ws.Cells.Item($formulaRow, 6)
    $fCell.Formula = $fFormula
    $fCell.Style = "Standaard"
}

Add-IncreaseBlock 18 19 "E3" "E10" "=((E3 / 114.202998) * 100) - 100" "=((E10 / 113.658804) * 100) - 100"
Add-IncreaseBlock 50 51 "E35" "E42" "=((E35 / 114.202998) * 100) - 100" "=((E42 / 113.658804) * 100) - 100"
Add-IncreaseBlock 82 83 "E67" "E74" "=((E67 / 114.202998) * 100) - 100" "=((E74 / 113.658804) * 100) - 100"
Add-IncreaseBlock 113 114 "" "" "=(D19 + D51 + D83) / 3" "=(F19 + F51 + F83) / 3"

# --- 3. Misc cosmetic bits captured by the diff: selection / scroll position.
$w = $excel.ActiveWindow
$ws.Range("H115").Select()
$w.ScrollRow = 87
$w.ScrollColumn = 1
$w.Top = 420
$w.Left = 800
